$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.917.59'
$ws.Range("E2").Value = '  +1.87%  '

$ws.Range("D3").Value = '1.812.16'
$ws.Range("E3").Value = '  +1.05%  '

$c = $ws.Range("D4")
$c.Value = "'1.002"
$c.Style = "Normal"
$ws.Range("E4").Value = '  -0.10%  '

$c = $ws.Range("D5")
$c.Value = "'338.09"
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.45%  '

$c = $ws.Range("D6")
$c.Value = "'0.9979"
$c.Style = "Normal"
$ws.Range("E6").Value = '  -0.22%  '

$c = $ws.Range("D7")
$c.Value = "'0.3928"
$c.Style = "Normal"
$ws.Range("E7").Value = '  +3.51%  '

$c = $ws.Range("D8")
$c.Value = "'0.3485"
$c.Style = "Normal"
$ws.Range("E8").Value = '  +1.63%  '

$c = $ws.Range("D9")
$c.Value = "'48.39"
$c.Style = "Normal"
$ws.Range("E9").Value = '  +0.80%  '

$c = $ws.Range("D10")
$c.Value = "'1.202"
$c.Style = "Normal"
$ws.Range("E10").Value = '  +0.04%  '

$c = $ws.Range("D11")
$c.Value = "'0.07575"
$c.Style = "Normal"
$ws.Range("E11").Value = '  +1.07%  '

$c = $ws.Range("D12")
$c.Value = "'0.9976"
$c.Style = "Normal"
$ws.Range("E12").Value = '  -0.25%  '

$c = $ws.Range("D13")
$c.Value = "'22.22"
$c.Style = "Normal"
$ws.Range("E13").Value = '  +0.76%  '

$c = $ws.Range("D14")
$c.Value = "'6.532"
$c.Style = "Normal"
$ws.Range("E14").Value = '  +0.86%  '

$ws.Range("D15").Value = '1.809.46'
$ws.Range("E15").Value = '  +1.10%  '

$c = $ws.Range("D16")
$c.Value = "'7.182"
$c.Style = "Normal"
$ws.Range("E16").Value = '  +1.87%  '

$c = $ws.Range("D17")
$c.Value = "'0.00001107"
$c.Style = "Normal"
$ws.Range("E17").Value = '  +0.89%  '

$c = $ws.Range("D18")
$c.Value = "'0.06685"
$c.Style = "Normal"
$ws.Range("E18").Value = '  +0.59%  '

$c = $ws.Range("D19")
$c.Value = "'85.14"
$c.Style = "Normal"
$ws.Range("E19").Value = '  +0.63%  '

$c = $ws.Range("D20")
$c.Value = "'0.9981"
$c.Style = "Normal"
$ws.Range("E20").Value = '  -0.22%  '

$ws.Range("E21").Value = '  +2.66%  '

$ws.Range("E22").Value = '  +1.30%  '

$ws.Range("D23").Value = '27.931.34'
$ws.Range("E23").Value = '  +2.03%  '

$ws.Range("E24").Value = '  +2.73%  '

$c = $ws.Range("D25")
$c.Value = "'2.401"
$c.Style = "Normal"
$ws.Range("E25").Value = '  -1.94%  '

$c = $ws.Range("D26")
$c.Value = "'2.552"
$c.Style = "Normal"
$ws.Range("E26").Value = '  -1.06%  '

$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Range("D27")
$c.Value = "'21.38"
$c.Style = "Normal"
$ws.Range("E27").Value = '  +0.02%  '

$ws.Range("B28").Value = 'ImmutableX'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range("D28")
$c.Value = "'1.475"
$c.Style = "Normal"
$ws.Range("E28").Value = '  -1.37%  '

$c = $ws.Range("D29")
$c.Value = "'154.79"
$c.Style = "Normal"
$ws.Range("E29").Value = '  +2.56%  '

$ws.Range("D30").Value = '2.015.06'
$ws.Range("E30").Value = '  +1.21%  '

$c = $ws.Range("D31")
$c.Value = "'135.41"
$c.Style = "Normal"
$ws.Range("E31").Value = '  +1.50%  '

$c = $ws.Range("D32")
$c.Value = "'4.040"
$c.Style = "Normal"
$ws.Range("E32").Value = '  -0.48%  '

$c = $ws.Range("D33")
$c.Value = "'6.112"
$c.Style = "Normal"
$ws.Range("E33").Value = '  -0.29%  '

$c = $ws.Range("D34")
$c.Value = "'0.08841"
$c.Style = "Normal"
$ws.Range("E34").Value = '  +1.54%  '

$ws.Range("E35").Value = '  -0.21%  '

$c = $ws.Range("D36")
$c.Value = "'5.513"
$c.Style = "Normal"
$ws.Range("E36").Value = '  +1.37%  '

$c = $ws.Range("D37")
$c.Value = "'0.6929"
$c.Style = "Normal"
$ws.Range("E37").Value = '  +0.46%  '

$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range("D38")
$c.Value = "'0.02427"
$c.Style = "Normal"
$ws.Range("E38").Value = '  +3.80%  '

$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range("D39")
$c.Value = "'0.06539"
$c.Style = "Normal"
$ws.Range("E39").Value = '  +3.01%  '

$ws.Range("B40").Value = 'WEMIXTOKEN'
$ws.Range("C40").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$c = $ws.Range("D40")
$c.Value = "'1.611"
$c.Style = "Normal"
$ws.Range("E40").Value = '  -3.94%  '

$c = $ws.Range("D41")
$c.Value = "'0.2229"
$c.Style = "Normal"
$ws.Range("E41").Value = '  +1.22%  '

$c = $ws.Range("D42")
$c.Value = "'1.263"
$c.Style = "Normal"
$ws.Range("E42").Value = '  -0.78%  '

$c = $ws.Range("D43")
$c.Value = "'8.529"
$c.Style = "Normal"
$ws.Range("E43").Value = '  -3.79%  '

$c = $ws.Range("D44")
$c.Value = "'14.76"
$c.Style = "Normal"
$ws.Range("E44").Value = '  +2.76%  '

$c = $ws.Range("D45")
$c.Value = "'0.6526"
$c.Style = "Normal"
$ws.Range("E45").Value = '  +1.11%  '

$ws.Range("E46").Value = '  -0.22%  '

$c = $ws.Range("D47")
$c.Value = "'3.872"
$c.Style = "Normal"
$ws.Range("E47").Value = '  +0.46%  '

$c = $ws.Range("D48")
$c.Value = "'2.166"
$c.Style = "Normal"
$ws.Range("E48").Value = '  +2.32%  '

$c = $ws.Range("D49")
$c.Value = "'132.76"
$c.Style = "Normal"
$ws.Range("E49").Value = '  +2.15%  '

$c = $ws.Range("D50")
$c.Value = "'0.07204"
$c.Style = "Normal"
$ws.Range("E50").Value = '  +0.35%  '

$c = $ws.Range("D51")
$c.Value = "'80.45"
$c.Style = "Normal"
$ws.Range("E51").Value = '  +1.67%  '
